# CU-1vc90uv api test for checking registering of user
# Adds a new "apiTest" worksheet (after the existing "LoginTest" sheet) that
# holds a single row of registration test data: email/password/firstName/lastName
# headers plus one data row, with the email cell carrying a mailto: hyperlink.

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

# Insert the new sheet right after LoginTest and make it the active tab.
$ws2 = $wb.Worksheets.Add($null, $ws1)
$ws2.Name = "apiTest"

# Header row
$ws2.Cells.Item(1, 1).Value = "email"
$ws2.Cells.Item(1, 2).Value = "password"
$ws2.Cells.Item(1, 3).Value = "firstName"
$ws2.Cells.Item(1, 4).Value = "lastName"

# Data row
$ws2.Cells.Item(2, 1).Value = "srdjan.rados@htecgroup.com"
$ws2.Cells.Item(2, 2).Value = "Qwertysha1@"
$ws2.Cells.Item(2, 3).Value = "Srdjan"
$ws2.Cells.Item(2, 4).Value = "Rados"

# Hyperlink the email address to a mailto: link, same as the LoginTest sheet does.
$ws2.Hyperlinks.Add($ws2.Range("A2"), "mailto:srdjan.rados@htecgroup.com", [Type]::Missing, [Type]::Missing, "srdjan.rados@htecgroup.com")

# Adding the hyperlink auto-applies the built-in "Hyperlink" font style; the
# source workbook keeps the plain default cell style, so clear that back off.
$ws2.Range("A2").ClearFormats()
